# Resolving version conflicts in Excel documents (I hate merge conflicts).
# Adds a new "url" Item_Attribute_Name block (mirrors the existing
# "description" block) plus its content_location/href row to the
# FutureShop Product_Tags sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "url" attribute block - same tag path as the existing "description"
# block (div.prod-info > h4.prod-title > a), ending in a content_location
# row whose value is now "href" (previously blank for "description").

$ws.Range("A46").Value = "Item_Attribute_Name"
$ws.Range("B46").Value = "url"

$ws.Range("A47").Value = "html_tag"
$ws.Range("B47").Value = "div"

$ws.Range("A48").Value = "html_tag_attribute_name"
$ws.Range("B48").Value = "class"

$ws.Range("A49").Value = "html_tag_attribute_value"
$ws.Range("B49").Value = "prod-info"

$ws.Range("A51").Value = "Item_Attribute_Name"
$ws.Range("B51").Value = "url"

$ws.Range("A52").Value = "html_tag"
$ws.Range("B52").Value = "h4"

$ws.Range("A53").Value = "html_tag_attribute_name"
$ws.Range("B53").Value = "class"

$ws.Range("A54").Value = "html_tag_attribute_value"
$ws.Range("B54").Value = "prod-title"

$ws.Range("A56").Value = "Item_Attribute_Name"
$ws.Range("B56").Value = "url"

$ws.Range("A57").Value = "html_tag"
$ws.Range("B57").Value = "a"

$ws.Range("A58").Value = "html_tag_attribute_name"

$ws.Range("A59").Value = "html_tag_attribute_value"

$ws.Range("A61").Value = "Item_Attribute_Name"
$ws.Range("B61").Value = "url"

$ws.Range("A62").Value = "content_location"
$ws.Range("B62").Value = "href"

# Match the saved view state from the diff: B62 selected (new bottom of data).
$ws.Range("B62").Select()
